$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.342.10"
$ws.Range("E2").Value = "  -2.00%  "

# Row 3
$ws.Range("D3").Value = "1.851.21"
$ws.Range("E3").Value = "  -1.44%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6971"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.43%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3070"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.76%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07478"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.73"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.99%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08129"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.29%  "

# Row 12
$ws.Range("D12").Value = "1.897.04"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7264"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.40%  "

# Row 14
$ws.Range("E14").Value = "  -4.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.58"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.50%  "

# Row 16
$ws.Range("D16").Value = "29.860.91"
$ws.Range("E16").Value = "  -0.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.918"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007742"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.51%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.195.61"
$ws.Range("E21").Value = "  +2.96%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.647"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1482"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.030"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.95%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.942"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.85%  "

# Row 30
$ws.Range("E30").Value = "  -7.62%  "

# Row 31
$ws.Range("E31").Value = "  -1.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.410"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.070"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05272"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.02%  "

# Row 35
$ws.Range("E35").Value = "  -3.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7213"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.20%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.682"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01869"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.78%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.719"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8872"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4311"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.43%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.920"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.70%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.044.71"
$ws.Range("E45").Value = "  -5.81%  "

# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.75%  "

# Row 48
$ws.Range("D48").Value = "2.037.28"
$ws.Range("E48").Value = "  +0.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.255"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.99%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.763"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.213"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.26%  "
